$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13: DAMSLTag "b" -> "ba", DialogAct "Acknowledge (Backchannel)" -> "Appreciation"
$ws.Range("I13").Value = "ba"
$ws.Range("J13").Value = "Appreciation"

# Row 24: DAMSLTag "sd" -> "sv", DialogAct "Statement-non-opinion" -> "Statement-opinion"
$ws.Range("I24").Value = "sv"
$ws.Range("J24").Value = "Statement-opinion"
